# Update LR-pair table (Edn1-Ednrb) to include the "ECs" sending/target cluster
# that was missing from the previous NATMI run (per Dr Hou's advice).
# The table becomes the full 3x3 cross-product of clusters: ECs, FAPs, sCs
# as both Sending cluster (col A) and Target cluster (col D), rows 2-10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Edn1"
$ws.Cells.Item(2, 3).Value = "Ednrb"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2.0
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 7.803861666666666
$ws.Cells.Item(2, 8).Value = 23.411585
$ws.Cells.Item(2, 9).Value = 0.7185921299200277
$ws.Cells.Item(2, 10).Value = 0.7185921299200276
$ws.Cells.Item(2, 11).Value = 3.0
$ws.Cells.Item(2, 12).Value = 1.0
$ws.Cells.Item(2, 13).Value = 26.46535033333333
$ws.Cells.Item(2, 14).Value = 79.396051
$ws.Cells.Item(2, 15).Value = 0.3009886585509795
$ws.Cells.Item(2, 16).Value = 0.3009886585509795
$ws.Cells.Item(2, 17).Value = 206.5319329612039
$ws.Cells.Item(2, 18).Value = 1858.787396650835
$ws.Cells.Item(2, 19).Value = 0.2162880812299203
$ws.Cells.Item(2, 20).Value = 0.2162880812299203

# Row 3: ECs -> FAPs
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Edn1"
$ws.Cells.Item(3, 3).Value = "Ednrb"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2.0
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 7.803861666666666
$ws.Cells.Item(3, 8).Value = 23.411585
$ws.Cells.Item(3, 9).Value = 0.7185921299200277
$ws.Cells.Item(3, 10).Value = 0.7185921299200276
$ws.Cells.Item(3, 11).Value = 1.0
$ws.Cells.Item(3, 12).Value = 0.3333333333333333
$ws.Cells.Item(3, 13).Value = 0.04671833333333333
$ws.Cells.Item(3, 14).Value = 0.140155
$ws.Cells.Item(3, 15).Value = 0.000531324479087915
$ws.Cells.Item(3, 16).Value = 0.000531324479087915
$ws.Cells.Item(3, 17).Value = 0.3645834106305555
$ws.Cells.Item(3, 18).Value = 3.281250695675
$ws.Cells.Item(3, 19).Value = 0.0003818055891064341
$ws.Cells.Item(3, 20).Value = 0.000381805589106434

# Row 4: ECs -> sCs
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Edn1"
$ws.Cells.Item(4, 3).Value = "Ednrb"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 2.0
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 7.803861666666666
$ws.Cells.Item(4, 8).Value = 23.411585
$ws.Cells.Item(4, 9).Value = 0.7185921299200277
$ws.Cells.Item(4, 10).Value = 0.7185921299200276
$ws.Cells.Item(4, 11).Value = 3.0
$ws.Cells.Item(4, 12).Value = 1.0
$ws.Cells.Item(4, 13).Value = 61.41599633333334
$ws.Cells.Item(4, 14).Value = 184.247989
$ws.Cells.Item(4, 15).Value = 0.6984800169699326
$ws.Cells.Item(4, 16).Value = 0.6984800169699326
$ws.Cells.Item(4, 17).Value = 479.2819395058406
$ws.Cells.Item(4, 18).Value = 4313.537455552565
$ws.Cells.Item(4, 19).Value = 0.501922243101001
$ws.Cells.Item(4, 20).Value = 0.5019222431010009

# Row 5: FAPs -> ECs
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Edn1"
$ws.Cells.Item(5, 3).Value = "Ednrb"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3.0
$ws.Cells.Item(5, 6).Value = 1.0
$ws.Cells.Item(5, 7).Value = 2.696109
$ws.Cells.Item(5, 8).Value = 8.088327
$ws.Cells.Item(5, 9).Value = 0.2482620517329206
$ws.Cells.Item(5, 10).Value = 0.2482620517329206
$ws.Cells.Item(5, 11).Value = 3.0
$ws.Cells.Item(5, 12).Value = 1.0
$ws.Cells.Item(5, 13).Value = 26.46535033333333
$ws.Cells.Item(5, 14).Value = 79.396051
$ws.Cells.Item(5, 15).Value = 0.3009886585509795
$ws.Cells.Item(5, 16).Value = 0.3009886585509795
$ws.Cells.Item(5, 17).Value = 71.35346922185299
$ws.Cells.Item(5, 18).Value = 642.181222996677
$ws.Cells.Item(5, 19).Value = 0.07472406192020566
$ws.Cells.Item(5, 20).Value = 0.07472406192020564

# Row 6: FAPs -> FAPs
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Edn1"
$ws.Cells.Item(6, 3).Value = "Ednrb"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3.0
$ws.Cells.Item(6, 6).Value = 1.0
$ws.Cells.Item(6, 7).Value = 2.696109
$ws.Cells.Item(6, 8).Value = 8.088327
$ws.Cells.Item(6, 9).Value = 0.2482620517329206
$ws.Cells.Item(6, 10).Value = 0.2482620517329206
$ws.Cells.Item(6, 11).Value = 1.0
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.04671833333333333
$ws.Cells.Item(6, 14).Value = 0.140155
$ws.Cells.Item(6, 15).Value = 0.000531324479087915
$ws.Cells.Item(6, 16).Value = 0.000531324479087915
$ws.Cells.Item(6, 17).Value = 0.125957718965
$ws.Cells.Item(6, 18).Value = 1.133619470685
$ws.Cells.Item(6, 19).Value = 0.0001319077053142911
$ws.Cells.Item(6, 20).Value = 0.000131907705314291

# Row 7: FAPs -> sCs
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Edn1"
$ws.Cells.Item(7, 3).Value = "Ednrb"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3.0
$ws.Cells.Item(7, 6).Value = 1.0
$ws.Cells.Item(7, 7).Value = 2.696109
$ws.Cells.Item(7, 8).Value = 8.088327
$ws.Cells.Item(7, 9).Value = 0.2482620517329206
$ws.Cells.Item(7, 10).Value = 0.2482620517329206
$ws.Cells.Item(7, 11).Value = 3.0
$ws.Cells.Item(7, 12).Value = 1.0
$ws.Cells.Item(7, 13).Value = 61.41599633333334
$ws.Cells.Item(7, 14).Value = 184.247989
$ws.Cells.Item(7, 15).Value = 0.6984800169699326
$ws.Cells.Item(7, 16).Value = 0.6984800169699326
$ws.Cells.Item(7, 17).Value = 165.584220458267
$ws.Cells.Item(7, 18).Value = 1490.257984124403
$ws.Cells.Item(7, 19).Value = 0.1734060821074007
$ws.Cells.Item(7, 20).Value = 0.1734060821074007

# Row 8: sCs -> ECs
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Edn1"
$ws.Cells.Item(8, 3).Value = "Ednrb"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3.0
$ws.Cells.Item(8, 6).Value = 1.0
$ws.Cells.Item(8, 7).Value = 0.3599613333333334
$ws.Cells.Item(8, 8).Value = 1.079884
$ws.Cells.Item(8, 9).Value = 0.03314581834705166
$ws.Cells.Item(8, 10).Value = 0.03314581834705165
$ws.Cells.Item(8, 11).Value = 3.0
$ws.Cells.Item(8, 12).Value = 1.0
$ws.Cells.Item(8, 13).Value = 26.46535033333333
$ws.Cells.Item(8, 14).Value = 79.396051
$ws.Cells.Item(8, 15).Value = 0.3009886585509795
$ws.Cells.Item(8, 16).Value = 0.3009886585509795
$ws.Cells.Item(8, 17).Value = 9.526502793120445
$ws.Cells.Item(8, 18).Value = 85.738525138084
$ws.Cells.Item(8, 19).Value = 0.009976515400853524
$ws.Cells.Item(8, 20).Value = 0.009976515400853522

# Row 9: sCs -> FAPs
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Edn1"
$ws.Cells.Item(9, 3).Value = "Ednrb"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3.0
$ws.Cells.Item(9, 6).Value = 1.0
$ws.Cells.Item(9, 7).Value = 0.3599613333333334
$ws.Cells.Item(9, 8).Value = 1.079884
$ws.Cells.Item(9, 9).Value = 0.03314581834705166
$ws.Cells.Item(9, 10).Value = 0.03314581834705165
$ws.Cells.Item(9, 11).Value = 1.0
$ws.Cells.Item(9, 12).Value = 0.3333333333333333
$ws.Cells.Item(9, 13).Value = 0.04671833333333333
$ws.Cells.Item(9, 14).Value = 0.140155
$ws.Cells.Item(9, 15).Value = 0.000531324479087915
$ws.Cells.Item(9, 16).Value = 0.000531324479087915
$ws.Cells.Item(9, 17).Value = 0.01681679355777778
$ws.Cells.Item(9, 18).Value = 0.15135114202
$ws.Cells.Item(9, 19).Value = 0.00001761118466718988
$ws.Cells.Item(9, 20).Value = 0.00001761118466718987

# Row 10: sCs -> sCs
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Edn1"
$ws.Cells.Item(10, 3).Value = "Ednrb"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3.0
$ws.Cells.Item(10, 6).Value = 1.0
$ws.Cells.Item(10, 7).Value = 0.3599613333333334
$ws.Cells.Item(10, 8).Value = 1.079884
$ws.Cells.Item(10, 9).Value = 0.03314581834705166
$ws.Cells.Item(10, 10).Value = 0.03314581834705165
$ws.Cells.Item(10, 11).Value = 3.0
$ws.Cells.Item(10, 12).Value = 1.0
$ws.Cells.Item(10, 13).Value = 61.41599633333334
$ws.Cells.Item(10, 14).Value = 184.247989
$ws.Cells.Item(10, 15).Value = 0.6984800169699326
$ws.Cells.Item(10, 16).Value = 0.6984800169699326
$ws.Cells.Item(10, 17).Value = 22.10738392814178
$ws.Cells.Item(10, 18).Value = 198.966455353276
$ws.Cells.Item(10, 19).Value = 0.02315169176153094
$ws.Cells.Item(10, 20).Value = 0.02315169176153094

